$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("all_tools")
$ws.Range("G10").Value = 813
$ws.Range("I10").Value = -0.07157910106056362
$ws.Range("J10").Value = 0.4916600219100405
$ws.Range("K10").Value = -0.08830673038161191
$ws.Range("L10").Value = 0.5419792899048017
$ws.Range("G11").Value = 813
$ws.Range("I11").Value = -0.05507056613029693
$ws.Range("J11").Value = 0.5800104076897017
$ws.Range("K11").Value = -0.07444201065306216
$ws.Range("L11").Value = 0.6074026216973724
$ws.Range("G12").Value = 813
$ws.Range("I12").Value = 0.1102028102074909
$ws.Range("J12").Value = 0.265108023071319
$ws.Range("K12").Value = 0.1660257192865707
$ws.Range("L12").Value = 0.2491959671429019
$ws.Range("G25").Value = 39
$ws.Range("I25").Value = -0.1807753815155468
$ws.Range("J25").Value = 0.3541954904764164
$ws.Range("K25").Value = -0.2576049186596542
$ws.Range("L25").Value = 0.3354345184285685
$ws.Range("G26").Value = 39
$ws.Range("I26").Value = -0.1807753815155468
$ws.Range("J26").Value = 0.3541954904764164
$ws.Range("K26").Value = -0.2666436877354316
$ws.Range("L26").Value = 0.3181414648703181
$ws.Range("G27").Value = 39
$ws.Range("I27").Value = 0.3539900381483285
$ws.Range("J27").Value = 0.0705613685189203
$ws.Range("K27").Value = 0.4341802833034056
$ws.Range("L27").Value = 0.09288178063084394
$ws.Range("G28").Value = 39
$ws.Range("K28").Value = -0.2493004677260264
$ws.Range("L28").Value = 0.3517858440384553
$ws.Range("G29").Value = 39
$ws.Range("K29").Value = -0.1491396897503261
$ws.Range("L29").Value = 0.5814513259975999

$ws = $wb.Worksheets.Item("infer")
$ws.Range("F10").Value = 23
$ws.Range("G10").Value = 24
$ws.Range("I10").Value = -0.1454025530693833
$ws.Range("J10").Value = 0.2372373518450496
$ws.Range("K10").Value = -0.17271903862684
$ws.Range("L10").Value = 0.2303502122764337
$ws.Range("F11").Value = 23
$ws.Range("G11").Value = 24
$ws.Range("I11").Value = -0.1395616700784287
$ws.Range("J11").Value = 0.2348980869048207
$ws.Range("K11").Value = -0.1674579385094694
$ws.Range("L11").Value = 0.2450782275649824
$ws.Range("F12").Value = 23
$ws.Range("G12").Value = 24
$ws.Range("I12").Value = 0.02140819589682411
$ws.Range("J12").Value = 0.8544862615484419
$ws.Range("K12").Value = 0.02708713119452734
$ws.Range("L12").Value = 0.8518765230635053
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("I25").Value = -0.3535533905932737
$ws.Range("J25").Value = 0.1037416782365415
$ws.Range("K25").Value = -0.4200840252084029
$ws.Range("L25").Value = 0.105228057983522
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("I26").Value = -0.1649915822768611
$ws.Range("J26").Value = 0.4476990724652935
$ws.Range("K26").Value = -0.1960392117639214
$ws.Range("L26").Value = 0.4668248490265503
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("I27").Value = 0.02366905341655754
$ws.Range("J27").Value = 0.9135633303377861
$ws.Range("K27").Value = 0.02802621677476181
$ws.Range("L27").Value = 0.9179387985999929
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("I28").Value = -0.2625754538144587
$ws.Range("J28").Value = 0.2314460271038938
$ws.Range("K28").Value = -0.3089716991054783
$ws.Range("L28").Value = 0.2442606266224961
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("I29").Value = 0.2592724864350675
$ws.Range("J29").Value = 0.2328233516916538
$ws.Range("K29").Value = 0.3080616184861621
$ws.Range("L29").Value = 0.2457251662216493
